$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate "Newburyport CR line (beyond Ipswich, closed) 1976" row (row 149)
$ws.Rows(149).Delete()

# Rename "Orange Line (Malden-Oak Grove)" to "Orange Line (Sullivan-Oak Grove)" (now row 149 after delete)
$ws.Range("B149").Value = "Orange Line (Sullivan–Oak Grove)"
